$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "67.205.40"
$ws.Range('E2').Value = "  -3.23%  "
$ws.Range('D3').Value = "3.700.05"
$ws.Range('E3').Value = "  -1.66%  "
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = "  +0.16%  "
$ws.Range('D5').Value = "'591.81"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "  -3.66%  "
$ws.Range('D6').Value = "'166.02"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "  -5.96%  "
$ws.Range('D7').Value = "3.696.98"
$ws.Range('E7').Value = "  -1.67%  "
$ws.Range('E8').Value = "  +0.07%  "
$ws.Range('D9').Value = "'0.523"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "  -0.64%  "
$ws.Range('D10').Value = "'0.161"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "  -2.74%  "
$ws.Range('D11').Value = "'6.15"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = "  -4.19%  "
$ws.Range('D12').Value = "'0.458"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "  -5.07%  "
$ws.Range('D13').Value = "'37.85"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "  -4.35%  "
$ws.Range('D14').Value = "'0.0000241"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "  -4.79%  "
$ws.Range('D15').Value = "4.340.03"
$ws.Range('E15').Value = "  -1.10%  "
$ws.Range('D16').Value = "3.715.20"
$ws.Range('E16').Value = "  -1.32%  "
$ws.Range('D17').Value = "67.375.73"
$ws.Range('E17').Value = "  -3.08%  "
$ws.Range('E18').Value = "  -3.62%  "
$ws.Range('D19').Value = "'7.09"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "  -5.75%  "
$ws.Range('D20').Value = "'17.12"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = "  +3.57%  "
$ws.Range('D21').Value = "'487.12"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "  -4.27%  "
$ws.Range('D22').Value = "'9.11"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "  -3.99%  "
$ws.Range('D23').Value = "'0.722"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = "  -1.37%  "
$ws.Range('D24').Value = "'85.02"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "  -1.42%  "
$ws.Range('B25').Value = "Fetch.AI"
$ws.Range('C25').Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range('D25').Value = "'2.32"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "  -5.95%  "
$ws.Range('B26').Value = "PEPE"
$ws.Range('C26').Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range('D26').Value = "'0.0000141"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "  +1.09%  "
$ws.Range('D27').Value = "'12.14"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "  -5.50%  "
$ws.Range('B28').Value = "Dai"
$ws.Range('C28').Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range('D28').Value = "'0.997"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "  -0.27%  "
$ws.Range('B29').Value = "RenderToken"
$ws.Range('C29').Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('D29').Value = "'9.98"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = "  -4.95%  "
$ws.Range('D30').Value = "'2.92"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "  -2.59%  "
$ws.Range('D31').Value = "'2.36"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "  -6.09%  "
$ws.Range('D32').Value = "'7.71"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = "  -5.31%  "
$ws.Range('D33').Value = "'31.64"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "  +2.01%  "
$ws.Range('D34').Value = "'0.107"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = "  -6.99%  "
$ws.Range('D35').Value = "'1.00"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = "  +0.25%  "
$ws.Range('D36').Value = "'0.992"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = "  -4.93%  "
$ws.Range('D37').Value = "'5.73"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "  -6.12%  "
$ws.Range('E38').Value = "  -6.21%  "
$ws.Range('D39').Value = "'0.322"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "  -5.11%  "
$ws.Range('D40').Value = "'444.87"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "  -5.52%  "
$ws.Range('D41').Value = "'48.97"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "  -1.58%  "
$ws.Range('D42').Value = "'1.96"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = "  -5.35%  "
$ws.Range('D43').Value = "'2.78"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = "  -7.06%  "
$ws.Range('D44').Value = "'8.28"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "  -3.32%  "
$ws.Range('D46').Value = "'39.76"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = "  -10.11%  "
$ws.Range('B47').Value = "Monero"
$ws.Range('C47').Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('D47').Value = "'140.71"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "  +0.91%  "
$ws.Range('B48').Value = "Maker"
$ws.Range('C48').Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('D48').Value = "2.785.53"
$ws.Range('E48').Value = "  -5.27%  "
$ws.Range('D49').Value = "'0.0345"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "  -4.49%  "
$ws.Range('B50').Value = "InjectiveProtocol"
$ws.Range('C50').Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range('D50').Value = "'25.34"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "  -8.08%  "
$ws.Range('B51').Value = "EnergySwap"
$ws.Range('C51').Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('D51').Value = "'23.75"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = "  +7.94%  "
